$d = $word.ActiveDocument

# The abstract originally says the integrated dataset draws on "11" boat-based
# surveys and spans "1959 - 2020". Update the survey count to 15 and extend
# the year range through 2021.

# "... water quality data from 11 boat-based ..." -> "... from 15 boat-based ..."
$found1 = $d.Content.Find.Execute(
    "water quality data from 11 boat-based", $true, $false, $false, $false, $false,
    $true, 1, $false, "water quality data from 15 boat-based", 2)

# "... (surface) from 1959 - 2020." -> "... (surface) from 1959 - 2021."
$found2 = $d.Content.Find.Execute(
    "from 1959 - 2020", $true, $false, $false, $false, $false,
    $true, 1, $false, "from 1959 - 2021", 2)

Write-Output ("survey count updated: " + $found1)
Write-Output ("year range updated: " + $found2)
